# Update cryptos list - Tue Feb 6 23:35:09 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# These D-column cells hold numeric-looking strings (e.g. "303.17") that must
# stay stored as plain text, exactly like the rest of the Price column.
# Force a Text number format first so Excel does not silently coerce them to
# real numbers when the new value is assigned.
$textCells = @(
    "D5","D6","D10","D11","D13","D14","D19","D22","D23","D24",
    "D25","D26","D27","D30","D33","D34","D36","D40","D43","D47",
    "D49","D50","D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
Set-Cell "D2" "43.102.56"
Set-Cell "E2" "  +1.11%  "

# Row 3 - Ethereum
Set-Cell "D3" "2.375.89"
Set-Cell "E3" "  +3.37%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  +0.06%  "

# Row 5 - BNB
Set-Cell "D5" "303.17"
Set-Cell "E5" "  +0.81%  "

# Row 6 - Solana
Set-Cell "D6" "96.93"
Set-Cell "E6" "  +1.38%  "

# Row 7 - XRP
Set-Cell "E7" "  -0.08%  "

# Row 8 - USDC
Set-Cell "E8" "  -0.08%  "

# Row 9 - Cardano
Set-Cell "E9" "  +1.27%  "

# Row 10 - Avalanche
Set-Cell "D10" "34.18"
Set-Cell "E10" "  -0.99%  "

# Row 11 - Dogecoin
Set-Cell "D11" "0.0785"
Set-Cell "E11" "  +0.45%  "

# Row 12 - TRON
Set-Cell "E12" "  +1.76%  "

# Row 13 - Chainlink
Set-Cell "D13" "18.31"
Set-Cell "E13" "  -4.34%  "

# Row 14 - Polkadot
Set-Cell "D14" "6.81"
Set-Cell "E14" "  +1.35%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Cell "D15" "2.748.08"
Set-Cell "E15" "  +3.59%  "

# Row 16 - WrappedEther
Set-Cell "D16" "2.370.43"
Set-Cell "E16" "  +3.11%  "

# Row 17 - Polygon
Set-Cell "E17" "  +3.61%  "

# Row 18 - WrappedBTC
Set-Cell "D18" "43.100.08"
Set-Cell "E18" "  +1.24%  "

# Row 19 - InternetComputer(DFINITY)
Set-Cell "D19" "12.16"
Set-Cell "E19" "  -0.90%  "

# Row 20 - Uniswap
Set-Cell "E20" "  +4.89%  "

# Row 21 - ShibaInu
Set-Cell "D21" "0.0₃0887"
Set-Cell "E21" "  -0.23%  "

# Row 22 - Litecoin
Set-Cell "D22" "68.37"
Set-Cell "E22" "  +1.10%  "

# Row 23 - BitcoinCash
Set-Cell "D23" "235.20"
Set-Cell "E23" "  -0.11%  "

# Row 24 - ImmutableX
Set-Cell "D24" "2.20"
Set-Cell "E24" "  -3.13%  "

# Row 25 & 26 swapped: Dai and PancakeSwap exchange places
Set-Cell "B25" "PancakeSwap"
Set-Cell "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-Cell "D25" "2.44"
Set-Cell "E25" "  +1.38%  "

Set-Cell "B26" "Dai"
Set-Cell "C26" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-Cell "D26" "1.00"
Set-Cell "E26" "  -0.18%  "

# Row 27 - EthereumClassic
Set-Cell "D27" "24.83"
Set-Cell "E27" "  +2.20%  "

# Row 28 - Toncoin
Set-Cell "E28" "  +0.41%  "

# Row 29 - Cosmos
Set-Cell "E29" "  +1.17%  "

# Row 30 - InjectiveProtocol
Set-Cell "D30" "31.30"
Set-Cell "E30" "  -1.47%  "

# Row 31 - FirstDigitalUSD
Set-Cell "E31" "  -0.01%  "

# Row 32 - Filecoin
Set-Cell "E32" "  +1.75%  "

# Row 33 - Hedera
Set-Cell "D33" "0.0749"
Set-Cell "E33" "  +8.12%  "

# Row 34 - Celestia
Set-Cell "D34" "17.40"
Set-Cell "E34" "  -1.41%  "

# Row 35 - Kaspa
Set-Cell "E35" "  +5.64%  "

# Row 36 - ARBITRUM
Set-Cell "D36" "1.83"
Set-Cell "E36" "  +5.90%  "

# Row 37 - WEMIXToken
Set-Cell "E37" "  -0.70%  "

# Row 38 - RenderToken
Set-Cell "E38" "  -2.77%  "

# Row 39 - LidoDAOToken
Set-Cell "E39" "  +3.83%  "

# Row 40 - EnergySwap
Set-Cell "D40" "22.27"
Set-Cell "E40" "  +10.88%  "

# Row 41 - Stellar
Set-Cell "E41" "  -0.28%  "

# Row 42 - Maker
Set-Cell "D42" "1.958.55"
Set-Cell "E42" "  +0.52%  "

# Row 43 - Monero
Set-Cell "D43" "103.98"
Set-Cell "E43" "  -36.93%  "

# Row 44 - VeChain
Set-Cell "E44" "  +1.10%  "

# Row 45 - ApeXProtocol
Set-Cell "E45" "  +2.08%  "

# Row 46 - NEARProtocol
Set-Cell "E46" "  +0.01%  "

# Row 47 - FraxShare
Set-Cell "D47" "9.14"
Set-Cell "E47" "  -10.98%  "

# Rows 48-51: new entry RocketPoolETH inserted at 48, pushing MultiversX, Stacks,
# BitcoinSV down a row; TrustWalletToken drops off the bottom.
Set-Cell "B48" "RocketPoolETH"
Set-Cell "C48" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-Cell "D48" "2.602.85"
Set-Cell "E48" "  +3.19%  "

Set-Cell "B49" "MultiversX"
Set-Cell "C49" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-Cell "D49" "52.57"
Set-Cell "E49" "  -0.96%  "

Set-Cell "B50" "Stacks"
Set-Cell "C50" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-Cell "D50" "1.50"
Set-Cell "E50" "  +1.78%  "

Set-Cell "B51" "BitcoinSV"
Set-Cell "C51" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-Cell "D51" "71.90"
Set-Cell "E51" "  +1.85%  "
